# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# worksheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4389
$wsExhibit.Range("F5").Value = 24
$wsExhibit.Range("F10").Value = 156
$wsExhibit.Range("F12").Value = 1645
$wsExhibit.Range("F14").Value = 3505

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4389
$wsAll.Range("F5").Value = 24
$wsAll.Range("F12").Value = 156
$wsAll.Range("F16").Value = 1645
$wsAll.Range("F18").Value = 3505
